# The codeforiati:group-code and codeforiati:group-name columns (C and D)
# were reordered in the source data, which - because the sheet only
# references string values positionally - results in every row's C and D
# values being swapped (including the header row, where the column
# headers themselves swap places).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value()
    $dVal = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}
